$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "24/10/2025"
$ws.Range("B11").Value = "Heerenveen"
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = "NAC Breda"
$ws.Range("F11").Value = "D"
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 3
$ws.Range("J11").Value = 2
$ws.Range("K11").Value = 1.38
$ws.Range("L11").Value = 1.81
$ws.Range("M11").Value = 11
$ws.Range("N11").Value = 10
$ws.Range("O11").Value = 4
$ws.Range("P11").Value = 4
